$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Replace-ParagraphInner($range, $innerXml) {
    # $range must be the paragraph's text range EXCLUDING the trailing paragraph mark.
    $range.Delete()
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="' + $wNs + '"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) "CREATE DATABASE ProjetoAvaliacaoProfessores2; " -> "...3; "
#    (the digit run "2" becomes "3", keep 3 separate runs)
# ---------------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("CREATE DATABASE ProjetoAvaliacaoProfessores2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pEnd = $d.Paragraphs(1).Range.End
    $pStart = $d.Paragraphs(1).Range.Start
    # locate enclosing paragraph via the found range
    $para = $full.Paragraphs(1)
    $pr = $para.Range
    $inner = $d.Range($pr.Start, $pr.End - 1)
    $xml1 = '<w:r xmlns:w="' + $wNs + '"><w:t>CREATE DATABASE ProjetoAvaliacaoProfessores</w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>3</w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t xml:space="preserve">; </w:t></w:r>'
    Replace-ParagraphInner $inner $xml1
}

# ---------------------------------------------------------------------------
# 2) "USE ProjetoAvaliacaoProfessores2;" -> "...3;" and relocate the
#    "_GoBack" bookmark to sit right before the final ";" run.
# ---------------------------------------------------------------------------
$full = $d.Content
$found2 = $full.Find.Execute("USE ProjetoAvaliacaoProfessores2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $full.Paragraphs(1)
    $pr2 = $para2.Range
    $inner2 = $d.Range($pr2.Start, $pr2.End - 1)
    $xml2 = '<w:r xmlns:w="' + $wNs + '"><w:t>USE ProjetoAvaliacaoProfessores</w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>3</w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>;</w:t></w:r>'
    Replace-ParagraphInner $inner2 $xml2
}

# Move the "_GoBack" bookmark: Bookmarks.Add re-registers the single
# "_GoBack" bookmark, removing any earlier one (and its matching end).
$full = $d.Content
$found3 = $full.Find.Execute("USE ProjetoAvaliacaoProfessores3;", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $semiPos = $full.End - 1
    $bmRange = $d.Range($semiPos, $semiPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 3) Drop the MD5() wrapper around the literal 'Carla Dias' string.
# ---------------------------------------------------------------------------
$full = $d.Content
$found4 = $full.Find.Execute("MD5('Carla Dias')", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $full.Text = "'Carla Dias'"
}

# ---------------------------------------------------------------------------
# 4) Expand "INSERT INTO `professores`(`ID`, `NOME`) VALUES (null,md5('Walter
#    Feitosa'));" into the spell-checked run/proofErr structure, dropping the
#    md5() wrapper around 'Walter Feitosa'.
# ---------------------------------------------------------------------------
$full = $d.Content
$found5 = $full.Find.Execute("INSERT INTO ``professores``(``ID``, ``NOME``) VALUES (null,md5('Walter Feitosa'));", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $target = $d.Range($full.Start, $full.End)
    $target.Delete()
    $xml5 = '<w:r xmlns:w="' + $wNs + '"><w:t>INSERT INTO `professores`(`ID`, `NOME`) VALUES (</w:t></w:r>' +
            '<w:proofErr xmlns:w="' + $wNs + '" w:type="spellStart"/>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>null</w:t></w:r>' +
            '<w:proofErr xmlns:w="' + $wNs + '" w:type="spellEnd"/>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>,</w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r xmlns:w="' + $wNs + '"><w:t>' + "'" + 'Walter Feitosa' + "'" + ');</w:t></w:r>'
    Replace-ParagraphInner $target $xml5
}

Write-Host "edit complete"
